$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $text) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextValue $ws 'D2' '27.220.45'
Set-TextValue $ws 'E2' '  +1.22%  '
Set-TextValue $ws 'D3' '1.651.21'
Set-TextValue $ws 'E3' '  +0.25%  '
Set-TextValue $ws 'E4' '  -0.24%  '
Set-TextValue $ws 'D5' '217.49'
Set-TextValue $ws 'E5' '  +0.15%  '
Set-TextValue $ws 'D6' '0.517'
Set-TextValue $ws 'E6' '  +2.53%  '
Set-TextValue $ws 'E7' '  -0.24%  '
Set-TextValue $ws 'D8' '0.257'
Set-TextValue $ws 'E8' '  +1.45%  '
Set-TextValue $ws 'D9' '0.0629'
Set-TextValue $ws 'E9' '  +1.37%  '
Set-TextValue $ws 'D10' '20.05'
Set-TextValue $ws 'E10' '  +1.43%  '
Set-TextValue $ws 'D11' '0.0850'
Set-TextValue $ws 'E11' '  +0.74%  '
Set-TextValue $ws 'D12' '1.882.17'
Set-TextValue $ws 'E12' '  +0.19%  '
Set-TextValue $ws 'D13' '1.651.69'
Set-TextValue $ws 'E13' '  +0.21%  '
Set-TextValue $ws 'D14' '4.15'
Set-TextValue $ws 'E14' '  +0.33%  '
Set-TextValue $ws 'D15' '0.543'
Set-TextValue $ws 'E15' '  +2.72%  '
Set-TextValue $ws 'D16' '67.71'
Set-TextValue $ws 'E16' '  +2.06%  '
Set-TextValue $ws 'D17' '27.229.58'
Set-TextValue $ws 'E17' '  +1.11%  '
Set-TextValue $ws 'E18' '  +1.24%  '
Set-TextValue $ws 'D19' '219.86'
Set-TextValue $ws 'E19' '  +0.51%  '
Set-TextValue $ws 'E20' '  -0.17%  '
Set-TextValue $ws 'E21' '  +3.01%  '
Set-TextValue $ws 'D22' '2.59'
Set-TextValue $ws 'E22' '  +6.59%  '
Set-TextValue $ws 'D23' '4.44'
Set-TextValue $ws 'E23' '  +0.94%  '
Set-TextValue $ws 'E24' '  +0.59%  '
Set-TextValue $ws 'D25' '147.62'
Set-TextValue $ws 'E25' '  +1.16%  '
Set-TextValue $ws 'E26' '  +2.73%  '
Set-TextValue $ws 'E27' '  -0.16%  '
Set-TextValue $ws 'E28' '  +0.16%  '
Set-TextValue $ws 'E29' '  -0.29%  '
Set-TextValue $ws 'E30' '  -0.41%  '
Set-TextValue $ws 'E31' '  -0.20%  '
Set-TextValue $ws 'E32' '  +0.66%  '
Set-TextValue $ws 'E33' '  +1.78%  '
Set-TextValue $ws 'E34' '  +1.46%  '
Set-TextValue $ws 'D35' '1.270.42'
Set-TextValue $ws 'E35' '  +1.91%  '
Set-TextValue $ws 'E36' '  +0.29%  '
Set-TextValue $ws 'E37' '  +1.55%  '
Set-TextValue $ws 'E38' '  +2.89%  '
Set-TextValue $ws 'D39' '0.848'
Set-TextValue $ws 'E39' '  +2.17%  '
Set-TextValue $ws 'E40' '  -0.18%  '
Set-TextValue $ws 'D41' '0.810'
Set-TextValue $ws 'E41' '  +0.21%  '
Set-TextValue $ws 'E43' '  +5.75%  '
Set-TextValue $ws 'D44' '1.792.04'
Set-TextValue $ws 'E44' '  +0.10%  '
Set-TextValue $ws 'D45' '62.16'
Set-TextValue $ws 'E45' '  +2.07%  '
Set-TextValue $ws 'D46' '91.79'
Set-TextValue $ws 'E46' '  +0.20%  '
Set-TextValue $ws 'E47' '  +0.70%  '
Set-TextValue $ws 'E48' '  -0.31%  '
Set-TextValue $ws 'E49' '  -0.14%  '
Set-TextValue $ws 'D50' '7.70'
Set-TextValue $ws 'E50' '  +1.19%  '
Set-TextValue $ws 'D51' '0.0974'
Set-TextValue $ws 'E51' '  +0.07%  '
